# "added the slow trebuchet to calculations"
# Adds a third projectile (slower trebuchet, V0 = C4 = 10) alongside the
# existing two (C2 = 25, C3 = 38), renames the existing "Xf" range columns
# to "Xf1"/"Xf2" so they no longer collide, and adds new X3/Y3 columns (N/O)
# with the same style as the existing X/Y columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update simulation inputs -------------------------------------------------
# Time step used throughout the sheet increases from 0.7 to 0.85
$ws.Range("F2").Value = 0.85

# New (third, "slow trebuchet") initial velocity
$ws.Range("C4").Value = 10

# Fix the Y1 formula so it references $C$2 absolutely, like the other formulas
$ws.Range("I2").Formula = '=$C$2*SIN($B2)*$F$2-$G$2*($F$2)^2/2+$E$2'

# --- Rename existing "Xf" headers so each range has a distinct label --------
$ws.Range("J1").Value = "Xf1"
$ws.Range("M1").Value = "Xf2"

# --- Add new headers for the slow-trebuchet columns --------------------------
$ws.Range("N1").Value = "X3"
$ws.Range("O1").Value = "Y3"

# --- Add the new X3 / Y3 formulas for rows 2-11 -------------------------------
# Row 2 holds its own (non-shared) formula, rows 3-11 share one formula group,
# mirroring the layout already used by the X1/Y1/X2/Y2 columns.
$ws.Range("N2").Formula = '=$C$4*COS($B2)*$F$2+$D$2'
$ws.Range("O2").Formula = '=$C$4*SIN($B2)*$F$2-$G$2*($F$2)^2/2+$E$2'
$ws.Range("N3:N11").Formula = '=$C$4*COS($B3)*$F$2+$D$2'
$ws.Range("O3:O11").Formula = '=$C$4*SIN($B3)*$F$2-$G$2*($F$2)^2/2+$E$2'

# Match the numeric formatting used by the other X/Y columns (2 decimals)
$ws.Range("N2:N11").NumberFormat = "0.00"
$ws.Range("O2:O11").NumberFormat = "0.00"

# --- Update the saved selection -----------------------------------------------
[void]$ws.Range("O2:O11").Select()

# --- Try to widen the workbook window (best effort) ---------------------------
$aw = $excel.ActiveWindow
$aw.Width = 27060
